$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7974165678998304
$ws.Cells.Item(2, 3).Value = 0.08833590539570935
$ws.Cells.Item(2, 5).Value = 0.1646067732740342
$ws.Cells.Item(2, 6).Value = 3.205780380381213
$ws.Cells.Item(2, 7).Value = 0.002551435093365505
$ws.Cells.Item(2, 10).Value = 0.1419999282367996
$ws.Cells.Item(2, 11).Value = 0.8078892130485826
$ws.Cells.Item(2, 13).Value = 0.4005836055683574
$ws.Cells.Item(2, 14).Value = 3.056378165865667

$ws.Cells.Item(3, 2).Value = 0.758356083015201
$ws.Cells.Item(3, 3).Value = 0.08216649578149315
$ws.Cells.Item(3, 5).Value = 0.161962656734314
$ws.Cells.Item(3, 6).Value = 3.174920944638288
$ws.Cells.Item(3, 7).Value = 0.002555540182250857
$ws.Cells.Item(3, 10).Value = 0.1423908577401889
$ws.Cells.Item(3, 11).Value = 0.7651800901005856
$ws.Cells.Item(3, 13).Value = 0.3865647714351752
$ws.Cells.Item(3, 14).Value = 3.064247697972036

$ws.Cells.Item(4, 2).Value = 0.7348552175842258
$ws.Cells.Item(4, 3).Value = 0.07843218227955617
$ws.Cells.Item(4, 5).Value = 0.1604331549223623
$ws.Cells.Item(4, 6).Value = 3.157491198159846
$ws.Cells.Item(4, 7).Value = 0.002558192972678227
$ws.Cells.Item(4, 10).Value = 0.1426675231180283
$ws.Cells.Item(4, 11).Value = 0.7394490462849319
$ws.Cells.Item(4, 13).Value = 0.378205603375271
$ws.Cells.Item(4, 14).Value = 3.069812667018795

$ws.Cells.Item(5, 2).Value = 0.7253994948787579
$ws.Cells.Item(5, 3).Value = 0.07692383711302853
$ws.Cells.Item(5, 5).Value = 0.1598335466216341
$ws.Cells.Item(5, 6).Value = 3.150769788022714
$ws.Cells.Item(5, 7).Value = 0.002559307369848392
$ws.Cells.Item(5, 10).Value = 0.1427894683138007
$ws.Cells.Item(5, 11).Value = 0.7290869425008566
$ws.Cells.Item(5, 13).Value = 0.37486162084911
$ws.Cells.Item(5, 14).Value = 3.072264479576361

$ws.Cells.Item(6, 2).Value = 0.7238366848620501
$ws.Cells.Item(6, 3).Value = 0.07667418527485381
$ws.Cells.Item(6, 5).Value = 0.1597354132896953
$ws.Cells.Item(6, 6).Value = 3.14967672496357
$ws.Cells.Item(6, 7).Value = 0.002559494433059656
$ws.Cells.Item(6, 10).Value = 0.1428102727032083
$ws.Cells.Item(6, 11).Value = 0.7273737725908518
$ws.Cells.Item(6, 13).Value = 0.37431012592846
$ws.Cells.Item(6, 14).Value = 3.072682708549436

$ws.Cells.Item(7, 2).Value = 0.7347272045050772
$ws.Cells.Item(7, 3).Value = 0.07841178600243381
$ws.Cells.Item(7, 5).Value = 0.1604249724774043
$ws.Cells.Item(7, 6).Value = 3.157399007298807
$ws.Cells.Item(7, 7).Value = 0.002558207866613258
$ws.Cells.Item(7, 10).Value = 0.1426691304690735
$ws.Cells.Item(7, 11).Value = 0.7393087997130294
$ws.Cells.Item(7, 13).Value = 0.3781602524510248
$ws.Cells.Item(7, 14).Value = 3.069844988240021

$ws.Cells.Item(8, 2).Value = 0.7838482239577615
$ws.Cells.Item(8, 3).Value = 0.08619748341266131
$ws.Cells.Item(8, 5).Value = 0.1636755954116396
$ws.Cells.Item(8, 6).Value = 3.194824723888274
$ws.Cells.Item(8, 7).Value = 0.002552823145081367
$ws.Cells.Item(8, 10).Value = 0.1421271114510994
$ws.Cells.Item(8, 11).Value = 0.7930606580086419
$ws.Cells.Item(8, 13).Value = 0.395698321410805
$ws.Cells.Item(8, 14).Value = 3.058939261869853

$ws.Cells.Item(9, 2).Value = 0.8840193265573077
$ws.Cells.Item(9, 3).Value = 0.1018966026306742
$ws.Cells.Item(9, 5).Value = 0.1707947048579044
$ws.Cells.Item(9, 6).Value = 3.28028657032516
$ws.Cells.Item(9, 7).Value = 0.002543308108356596
$ws.Cells.Item(9, 10).Value = 0.1413553749022611
$ws.Cells.Item(9, 11).Value = 0.902398138992254
$ws.Cells.Item(9, 13).Value = 0.4320655779745124
$ws.Cells.Item(9, 14).Value = 3.04338334812465

$ws.Cells.Item(10, 2).Value = 0.9599916247516944
$ws.Cells.Item(10, 3).Value = 0.1137020496832406
$ws.Cells.Item(10, 5).Value = 0.1764783610762919
$ws.Cells.Item(10, 6).Value = 3.350479487285099
$ws.Cells.Item(10, 7).Value = 0.002536947147707674
$ws.Cells.Item(10, 10).Value = 0.1409666341701765
$ws.Cells.Item(10, 11).Value = 0.9851647880634289
$ws.Cells.Item(10, 13).Value = 0.4599975476568687
$ws.Cells.Item(10, 14).Value = 3.035528975784899

$ws.Cells.Item(11, 2).Value = 0.9950771101892997
$ws.Cells.Item(11, 3).Value = 0.1191334506376052
$ws.Cells.Item(11, 5).Value = 0.1791623379024685
$ws.Cells.Item(11, 6).Value = 3.384030539546131
$ws.Cells.Item(11, 7).Value = 0.0025341886415407
$ws.Cells.Item(11, 10).Value = 0.140828664746973
$ws.Cells.Item(11, 11).Value = 1.023355915202785
$ws.Cells.Item(11, 13).Value = 0.4729700678933213
$ws.Cells.Item(11, 14).Value = 3.032736544231383

$ws.Cells.Item(12, 2).Value = 1.008439050703601
$ws.Cells.Item(12, 3).Value = 0.1211990862408072
$ws.Cells.Item(12, 5).Value = 0.1801928239326855
$ws.Cells.Item(12, 6).Value = 3.396969063053405
$ws.Cells.Item(12, 7).Value = 0.002533163385515413
$ws.Cells.Item(12, 10).Value = 0.140782022753708
$ws.Cells.Item(12, 11).Value = 1.037896145337783
$ws.Cells.Item(12, 13).Value = 0.4779207979928302
$ws.Cells.Item(12, 14).Value = 3.031791733152502

$ws.Cells.Item(13, 2).Value = 1.005557940540257
$ws.Cells.Item(13, 3).Value = 0.1207538179582457
$ws.Cells.Item(13, 5).Value = 0.1799702627901141
$ws.Cells.Item(13, 6).Value = 3.394172128510917
$ws.Cells.Item(13, 7).Value = 0.00253338333468738
$ws.Cells.Item(13, 10).Value = 0.1407918184670365
$ws.Cells.Item(13, 11).Value = 1.034761167503916
$ws.Cells.Item(13, 13).Value = 0.4768528646546031
$ws.Cells.Item(13, 14).Value = 3.031990200137912

$ws.Cells.Item(14, 2).Value = 0.9961748836698803
$ws.Cells.Item(14, 3).Value = 0.1193032131133975
$ws.Cells.Item(14, 5).Value = 0.1792468337763466
$ws.Cells.Item(14, 6).Value = 3.385090317579738
$ws.Cells.Item(14, 7).Value = 0.002534103906469687
$ws.Cells.Item(14, 10).Value = 0.1408247150849178
$ws.Cells.Item(14, 11).Value = 1.02455058151935
$ws.Cells.Item(14, 13).Value = 0.4733765991903311
$ws.Cells.Item(14, 14).Value = 3.03265655475569

$ws.Cells.Item(15, 2).Value = 0.9904373742701296
$ws.Cells.Item(15, 3).Value = 0.1184158355596026
$ws.Cells.Item(15, 5).Value = 0.1788055507193604
$ws.Cells.Item(15, 6).Value = 3.379557866298512
$ws.Cells.Item(15, 7).Value = 0.002534547791865763
$ws.Cells.Item(15, 10).Value = 0.1408455954591226
$ws.Cells.Item(15, 11).Value = 1.018306483127361
$ws.Cells.Item(15, 13).Value = 0.471252277899417
$ws.Cells.Item(15, 14).Value = 3.033079394240275

$ws.Cells.Item(16, 2).Value = 0.9577093009161217
$ws.Cells.Item(16, 3).Value = 0.1133483320670621
$ws.Cells.Item(16, 5).Value = 0.1763049349835555
$ws.Cells.Item(16, 6).Value = 3.348319488894361
$ws.Cells.Item(16, 7).Value = 0.002537130133183714
$ws.Cells.Item(16, 10).Value = 0.140976434204255
$ws.Cells.Item(16, 11).Value = 0.9826798119327123
$ws.Cells.Item(16, 13).Value = 0.4591551222734651
$ws.Cells.Item(16, 14).Value = 3.03572720687329

$ws.Cells.Item(17, 2).Value = 0.9377664357302535
$ws.Cells.Item(17, 3).Value = 0.1102552910942052
$ws.Cells.Item(17, 5).Value = 0.1747960791733902
$ws.Cells.Item(17, 6).Value = 3.329571030348944
$ws.Cells.Item(17, 7).Value = 0.002538748853406954
$ws.Cells.Item(17, 10).Value = 0.1410666654501114
$ws.Cells.Item(17, 11).Value = 0.9609626370184401
$ws.Cells.Item(17, 13).Value = 0.4518021049703762
$ws.Cells.Item(17, 14).Value = 3.0375517422858

$ws.Cells.Item(18, 2).Value = 0.9263452536414718
$ws.Cells.Item(18, 3).Value = 0.1084819951736051
$ws.Cells.Item(18, 5).Value = 0.1739374941530585
$ws.Cells.Item(18, 6).Value = 3.318939857609621
$ws.Cells.Item(18, 7).Value = 0.002539692623265959
$ws.Cells.Item(18, 10).Value = 0.1411222218847854
$ws.Cells.Item(18, 11).Value = 0.9485223104741181
$ws.Cells.Item(18, 13).Value = 0.4475978899150306
$ws.Cells.Item(18, 14).Value = 3.038674614059758

$ws.Cells.Item(19, 2).Value = 0.922486717946299
$ws.Cells.Item(19, 3).Value = 0.107882570256379
$ws.Cells.Item(19, 5).Value = 0.1736483850488213
$ws.Cells.Item(19, 6).Value = 3.315366491969513
$ws.Cells.Item(19, 7).Value = 0.002540014355977727
$ws.Cells.Item(19, 10).Value = 0.1411416601351014
$ws.Cells.Item(19, 11).Value = 0.9443189455234915
$ws.Cells.Item(19, 13).Value = 0.4461787151508858
$ws.Cells.Item(19, 14).Value = 3.039067402329479

$ws.Cells.Item(20, 2).Value = 0.9398842702434251
$ws.Cells.Item(20, 3).Value = 0.1105839558257173
$ws.Cells.Item(20, 5).Value = 0.1749557404077535
$ws.Cells.Item(20, 6).Value = 3.331551054607928
$ws.Cells.Item(20, 7).Value = 0.002538575221615642
$ws.Cells.Item(20, 10).Value = 0.1410566815072443
$ws.Cells.Item(20, 11).Value = 0.9632692049851528
$ws.Cells.Item(20, 13).Value = 0.4525822536989565
$ws.Cells.Item(20, 14).Value = 3.037349913182226

$ws.Cells.Item(21, 2).Value = 0.9989288546055377
$ws.Cells.Item(21, 3).Value = 0.1197290490305534
$ws.Cells.Item(21, 5).Value = 0.1794589394168256
$ws.Cells.Item(21, 6).Value = 3.387751525116045
$ws.Cells.Item(21, 7).Value = 0.002533891732984519
$ws.Cells.Item(21, 10).Value = 0.1408149003418693
$ws.Cells.Item(21, 11).Value = 1.02754755582356
$ws.Cells.Item(21, 13).Value = 0.474396622269154
$ws.Cells.Item(21, 14).Value = 3.032457770677155

$ws.Cells.Item(22, 2).Value = 1.03796003086012
$ws.Cells.Item(22, 3).Value = 0.1257577160875769
$ws.Cells.Item(22, 5).Value = 0.1824843397837483
$ws.Cells.Item(22, 6).Value = 3.425842863275335
$ws.Cells.Item(22, 7).Value = 0.002530943428829117
$ws.Cells.Item(22, 10).Value = 0.1406895539140969
$ws.Cells.Item(22, 11).Value = 1.070012472353568
$ws.Cells.Item(22, 13).Value = 0.4888769461203424
$ws.Cells.Item(22, 14).Value = 3.029917100068786

$ws.Cells.Item(23, 2).Value = 1.017087805095912
$ws.Cells.Item(23, 3).Value = 0.1225353267820992
$ws.Cells.Item(23, 5).Value = 0.1808621076998094
$ws.Cells.Item(23, 6).Value = 3.405388084548633
$ws.Cells.Item(23, 7).Value = 0.002532506721981844
$ws.Cells.Item(23, 10).Value = 0.1407534593476285
$ws.Cells.Item(23, 11).Value = 1.047306351661376
$ws.Cells.Item(23, 13).Value = 0.4811280690324011
$ws.Cells.Item(23, 14).Value = 3.031212896034134

$ws.Cells.Item(24, 2).Value = 0.9389266594918411
$ws.Cells.Item(24, 3).Value = 0.1104353511385909
$ws.Cells.Item(24, 5).Value = 0.1748835299135223
$ws.Cells.Item(24, 6).Value = 3.330655426024975
$ws.Cells.Item(24, 7).Value = 0.00253865367953515
$ws.Cells.Item(24, 10).Value = 0.141061183782238
$ws.Cells.Item(24, 11).Value = 0.9622262649525339
$ws.Cells.Item(24, 13).Value = 0.4522294764092649
$ws.Cells.Item(24, 14).Value = 3.037440929875999

$ws.Cells.Item(25, 2).Value = 0.856505416866014
$ws.Cells.Item(25, 3).Value = 0.0976025369011353
$ws.Cells.Item(25, 5).Value = 0.1687891534386736
$ws.Cells.Item(25, 6).Value = 3.255869989932521
$ws.Cells.Item(25, 7).Value = 0.002545771089895288
$ws.Cells.Item(25, 10).Value = 0.1415329024444176
$ws.Cells.Item(25, 11).Value = 0.8723945443021819
$ws.Cells.Item(25, 13).Value = 0.3781602524510248
$ws.Cells.Item(25, 14).Value = 3.046965391872547
